$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (values) and column C (hour) for rows 2-12
$ws.Range("B2").Value = 1501.5
$ws.Range("C2").Value = 12

$ws.Range("B3").Value = 148
$ws.Range("C3").Value = 12

$ws.Range("C4").Value = 12

$ws.Range("B5").Value = 926
$ws.Range("C5").Value = 12

$ws.Range("B6").Value = 642
$ws.Range("C6").Value = 12

$ws.Range("B7").Value = 196
$ws.Range("C7").Value = 12

$ws.Range("C8").Value = 12

$ws.Range("B9").Value = 788
$ws.Range("C9").Value = 12

$ws.Range("B10").Value = 95
$ws.Range("C10").Value = 12

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = 12

$ws.Range("B12").Value = 85
$ws.Range("C12").Value = 12
